$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("GLOBAL RESULTS")
$ws.Range("C6").Value = 22065.220353469886
$ws.Range("C7").Value = 21685.220353469886
$ws.Range("C8").Value = 21403.263742865784
$ws.Range("C12").Value = 3051.422573251716
$ws.Range("C14").Value = 19013.79778021817
$ws.Range("C15").Value = 18633.79778021817
$ws.Range("C16").Value = 12173.797780218172
$ws.Range("C17").Value = 11944.25413471817
$ws.Range("C18").Value = 11323.830134718173
$ws.Range("C23").Value = 216385.89317935536
$ws.Range("C24").Value = 212659.36617935536
$ws.Range("C25").Value = 209894.31638397468
$ws.Range("C30").Value = 186461.66000137647
$ws.Range("C31").Value = 182735.13300137647
$ws.Range("C32").Value = 119384.1740013765
$ws.Range("C33").Value = 117133.11981023391
$ws.Range("C34").Value = 111048.83879063395

$ws = $wb.Worksheets.Item("FUSELAGE")
$ws.Range("C6").Value = 3610.0
$ws.Range("D6").Value = 38.49326236033189
$ws.Range("C7").Value = 2405.0
$ws.Range("D7").Value = -7.735098067424319
$ws.Range("C8").Value = 2837.0
$ws.Range("D8").Value = 8.83805687431069
$ws.Range("C9").Value = 2578.0
$ws.Range("D9").Value = -1.098163333812846
$ws.Range("C12").Value = 2936.5
$ws.Range("D12").Value = 12.655253440751956

$ws = $wb.Worksheets.Item("WING")
$ws.Range("C7").Value = 2308.0
$ws.Range("D7").Value = 32.81542224140415
$ws.Range("C9").Value = 2012.0
$ws.Range("D9").Value = 15.781901884620947
$ws.Range("C13").Value = 1904.1428571428569
$ws.Range("D13").Value = 9.575189592453313

$ws = $wb.Worksheets.Item("HORIZONTAL TAIL")
$ws.Range("C9").Value = 123.0
$ws.Range("D9").Value = -52.81254495756005
$ws.Range("C10").Value = 176.66666666666663
$ws.Range("D10").Value = -32.22398056234912

$ws = $wb.Worksheets.Item("NACELLES")
$ws.Range("C10").Value = 553.0
$ws.Range("D10").Value = 324.3034575360861
$ws.Range("C17").Value = 553.0
$ws.Range("D17").Value = 324.3034575360861

$ws = $wb.Worksheets.Item("LANDING GEARS")
$ws.Range("C5").Value = 737.0
$ws.Range("D5").Value = 6.027909653287345
$ws.Range("C6").Value = 883.0
$ws.Range("D6").Value = 27.03208171486123
$ws.Range("C7").Value = 1003.0
$ws.Range("D7").Value = 44.295784779168535
$ws.Range("C8").Value = 894.0
$ws.Range("D8").Value = 28.6145878290894
$ws.Range("C9").Value = 879.25
$ws.Range("D9").Value = 26.49259099410159
